# Remove the hyperlink emails ("show in console" / not clickable links any more)
# and update the email addresses, plus add the newest feedback row (shan / WAS).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two stale/incorrect e-mail addresses to the corrected one.
$ws.Range("C2").Value = "ranaabobakar777@gmail.com"
$ws.Range("C3").Value = "ranaabobakar777@gmail.com"

# Append the new row of feedback data.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "shan"
$ws.Range("C5").Value = "shan@gmail.com"
$ws.Range("D5").Value = "WAS"
$ws.Rows.Item(5).RowHeight = 12.8

# Remove the mailto: hyperlinks that used to be attached to the email column.
$ws.Hyperlinks.Delete()

# Restore selection as left by the author.
$ws.Range("C7").Select() | Out-Null
